$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to be stored as text even when the value looks numeric
    # (mirrors what a user does in Excel: format the cell as Text, type the
    # value, then restore the cell style so no residual formatting is left).
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

function Set-BlankTextCell($range) {
    # Create a present-but-empty text cell (matches a blank cell that was
    # nonetheless typed/formatted as text, e.g. an empty data-entry field).
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = "'"
    $r.Style = "Normal"
}

# --- Row 5: convert C5:M5 from text to numeric values ---
$ws.Range("C5").Value = 100
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
# N5, O5, P5 stay "t" (unchanged)

# --- Row 6 (new) ---
Set-BlankTextCell "A6"
Set-BlankTextCell "B6"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
Set-TextValue "N6" "0"
Set-TextValue "O6" "0"
Set-TextValue "P6" "0"

# --- Row 7 (new) ---
Set-BlankTextCell "A7"
Set-BlankTextCell "B7"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
Set-TextValue "N7" "0"
Set-TextValue "O7" "0"
Set-TextValue "P7" "0"

# --- Row 8 (new) ---
Set-TextValue "A8" "2025-03-31 14:29:35"
Set-TextValue "B8" "SHIFT_1"
Set-TextValue "C8" "1"
Set-TextValue "D8" "0"
Set-TextValue "E8" "0"
Set-TextValue "F8" "0"
Set-TextValue "G8" "0"
Set-TextValue "H8" "0"
Set-TextValue "I8" "0"
Set-TextValue "J8" "0"
Set-TextValue "K8" "0"
Set-TextValue "L8" "0"
Set-TextValue "M8" "0"
Set-TextValue "N8" "0"
Set-TextValue "O8" "0"
Set-TextValue "P8" "0"
